$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 values (identifier / title / levelOfDescription / extentAndMedium / notes)
$ws.Range("A2").Value = "MCH139"
$ws.Range("C2").Value = "FANAGALO (KITCHEN KAFFIR), ZULU VOCABULARY AND PHRASE BOOK"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# Apply the "Calibri 10 theme1" font used for the new row to all touched cells,
# including the two cells that stay textually empty (D2 / H2).
$newRowCells = @("A2", "C2", "D2", "E2", "F2", "G2", "H2")
foreach ($addr in $newRowCells) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.Font.ThemeColor = 1
}

# F2 ends up with an explicit (default) alignment flag in the target file.
$ws.Range("F2").WrapText = $false

# Select A2:H2 (matches the new active selection) while keeping the existing
# frozen header row/pane intact.
$ws.Range("A2:H2").Select()
$excel.ActiveWindow.FreezePanes = $true
